$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 61659
$ws.Range("B2").Value = "Steven Nguyen"
$ws.Range("C2").Value = "iboyer@example.com"
$ws.Range("D2").Value = "HR"
$ws.Range("E2").Value = 69167.67

$ws.Range("A3").Value = 1067
$ws.Range("B3").Value = "Amanda Torres"
$ws.Range("C3").Value = "lisa77@example.org"
$ws.Range("D3").Value = "Marketing"
$ws.Range("E3").Value = 74755.39

$ws.Range("A4").Value = 77832
$ws.Range("B4").Value = "Susan Williams"
$ws.Range("C4").Value = "munozphillip@example.org"
$ws.Range("D4").Value = "IT"
$ws.Range("E4").Value = 72115.27

$ws.Range("A5").Value = 19325
$ws.Range("B5").Value = "Tiffany Mathews"
$ws.Range("C5").Value = "travisclark@example.com"
$ws.Range("D5").Value = "Marketing"
$ws.Range("E5").Value = 64673.01

$ws.Range("A6").Value = 53246
$ws.Range("B6").Value = "Amy Boone"
$ws.Range("C6").Value = "sotojoseph@example.net"
$ws.Range("D6").Value = "IT"
$ws.Range("E6").Value = 97488.24000000001

$ws.Range("A7").Value = 69060
$ws.Range("B7").Value = "Gregory Delgado"
$ws.Range("C7").Value = "erikrobles@example.org"
$ws.Range("D7").Value = "HR"
$ws.Range("E7").Value = 74642.62

$ws.Range("A8").Value = 21824
$ws.Range("B8").Value = "Troy Jackson"
$ws.Range("C8").Value = "alijames@example.org"
$ws.Range("D8").Value = "HR"
$ws.Range("E8").Value = 97886.82000000001

$ws.Range("A9").Value = 77354
$ws.Range("B9").Value = "Alexis Morris"
$ws.Range("C9").Value = "kimwillis@example.org"
$ws.Range("D9").Value = "Marketing"
$ws.Range("E9").Value = 38066.89

$ws.Range("A10").Value = 92184
$ws.Range("B10").Value = "Jessica Alexander"
$ws.Range("C10").Value = "stephenhunt@example.com"
$ws.Range("D10").Value = "HR"
$ws.Range("E10").Value = 78584.07000000001

$ws.Range("A11").Value = 80488
$ws.Range("B11").Value = "Nicole Barnett"
$ws.Range("C11").Value = "kathrynbrewer@example.net"
$ws.Range("D11").Value = "IT"
$ws.Range("E11").Value = 92242.69

$ws.Range("A12").Value = 59838
$ws.Range("B12").Value = "Danielle Stokes"
$ws.Range("C12").Value = "michael66@example.net"
$ws.Range("D12").Value = "Marketing"
$ws.Range("E12").Value = 33535.64

$ws.Range("A13").Value = 23895
$ws.Range("B13").Value = "Daniel Parker"
$ws.Range("C13").Value = "terry00@example.com"
$ws.Range("D13").Value = "IT"
$ws.Range("E13").Value = 64770.08

$ws.Range("A14").Value = 49451
$ws.Range("B14").Value = "James Cruz"
$ws.Range("C14").Value = "jacob86@example.com"
$ws.Range("D14").Value = "Finance"
$ws.Range("E14").Value = 91037.25

$ws.Range("A15").Value = 88147
$ws.Range("B15").Value = "John Lawrence"
$ws.Range("C15").Value = "williamsjessica@example.com"
$ws.Range("D15").Value = "IT"
$ws.Range("E15").Value = 62773.63

$ws.Range("A16").Value = 75155
$ws.Range("B16").Value = "Christopher Taylor"
$ws.Range("C16").Value = "candace16@example.com"
$ws.Range("D16").Value = "Marketing"
$ws.Range("E16").Value = 73522.67999999999

$ws.Range("A17").Value = 13759
$ws.Range("B17").Value = "Briana Callahan"
$ws.Range("C17").Value = "brendabrown@example.net"
$ws.Range("D17").Value = "Marketing"
$ws.Range("E17").Value = 67722.32000000001

$ws.Range("A18").Value = 82142
$ws.Range("B18").Value = "Paul Davis"
$ws.Range("C18").Value = "qhaas@example.com"
$ws.Range("D18").Value = "HR"
$ws.Range("E18").Value = 76547.41

$ws.Range("A19").Value = 95753
$ws.Range("B19").Value = "Pamela Schultz"
$ws.Range("C19").Value = "aswanson@example.org"
$ws.Range("D19").Value = "IT"
$ws.Range("E19").Value = 64147.81

$ws.Range("A20").Value = 71361
$ws.Range("B20").Value = "Aaron Snyder"
$ws.Range("C20").Value = "alexis71@example.net"
$ws.Range("D20").Value = "Marketing"
$ws.Range("E20").Value = 51497.23

$ws.Range("A21").Value = 17994
$ws.Range("B21").Value = "Courtney Jordan"
$ws.Range("C21").Value = "elizabethgillespie@example.org"
$ws.Range("D21").Value = "Marketing"
$ws.Range("E21").Value = 82466.96000000001

$ws.Range("A22").Value = 78798
$ws.Range("B22").Value = "Valerie Johnson"
$ws.Range("C22").Value = "wagnertracey@example.org"
$ws.Range("D22").Value = "Marketing"
$ws.Range("E22").Value = 50435.27

$ws.Range("A23").Value = 38804
$ws.Range("B23").Value = "Brian Phillips"
$ws.Range("C23").Value = "wmassey@example.net"
$ws.Range("D23").Value = "Marketing"
$ws.Range("E23").Value = 76496.41

$ws.Range("A24").Value = 50540
$ws.Range("B24").Value = "Paula Henderson"
$ws.Range("C24").Value = "smithcrystal@example.net"
$ws.Range("D24").Value = "Finance"
$ws.Range("E24").Value = 55160.56

$ws.Range("A25").Value = 45778
$ws.Range("B25").Value = "Christina Ponce"
$ws.Range("C25").Value = "thomasbautista@example.org"
$ws.Range("D25").Value = "Marketing"
$ws.Range("E25").Value = 35877.29

$ws.Range("A26").Value = 9855
$ws.Range("B26").Value = "Colton Ball"
$ws.Range("C26").Value = "ahicks@example.net"
$ws.Range("D26").Value = "HR"
$ws.Range("E26").Value = 48572.78

$ws.Range("A27").Value = 83925
$ws.Range("B27").Value = "Kevin Diaz"
$ws.Range("C27").Value = "umcpherson@example.com"
$ws.Range("D27").Value = "Marketing"
$ws.Range("E27").Value = 75219.74000000001

$ws.Range("A28").Value = 9983
$ws.Range("B28").Value = "Tamara Diaz"
$ws.Range("C28").Value = "mcfarlandkayla@example.com"
$ws.Range("D28").Value = "HR"
$ws.Range("E28").Value = 55072.3

$ws.Range("A29").Value = 96866
$ws.Range("B29").Value = "Timothy Malone"
$ws.Range("C29").Value = "dyoung@example.net"
$ws.Range("D29").Value = "HR"
$ws.Range("E29").Value = 87109.81

$ws.Range("A30").Value = 7240
$ws.Range("B30").Value = "Samuel Gilbert"
$ws.Range("C30").Value = "xgreen@example.com"
$ws.Range("D30").Value = "HR"
$ws.Range("E30").Value = 97037.83

$ws.Range("A31").Value = 41500
$ws.Range("B31").Value = "Scott Mosley"
$ws.Range("C31").Value = "belljeanette@example.org"
$ws.Range("D31").Value = "Finance"
$ws.Range("E31").Value = 66865.72

$ws.Range("A32").Value = 95194
$ws.Range("B32").Value = "Michael Sanchez"
$ws.Range("C32").Value = "brian22@example.org"
$ws.Range("D32").Value = "IT"
$ws.Range("E32").Value = 75863.83

$ws.Range("A33").Value = 5528
$ws.Range("B33").Value = "Natalie Hobbs"
$ws.Range("C33").Value = "tcampos@example.net"
$ws.Range("D33").Value = "Finance"
$ws.Range("E33").Value = 69777.77

$ws.Range("A34").Value = 81695
$ws.Range("B34").Value = "Robert Williamson"
$ws.Range("C34").Value = "joseph09@example.org"
$ws.Range("D34").Value = "Finance"
$ws.Range("E34").Value = 42403.21

$ws.Range("A35").Value = 11207
$ws.Range("B35").Value = "John Fox"
$ws.Range("C35").Value = "colemanmadison@example.com"
$ws.Range("D35").Value = "Finance"
$ws.Range("E35").Value = 93745.14

$ws.Range("A36").Value = 77409
$ws.Range("B36").Value = "Jason Robinson"
$ws.Range("C36").Value = "desireemartinez@example.net"
$ws.Range("D36").Value = "IT"
$ws.Range("E36").Value = 96146.98

$ws.Range("A37").Value = 17129
$ws.Range("B37").Value = "Courtney Atkins"
$ws.Range("C37").Value = "zwalker@example.com"
$ws.Range("D37").Value = "IT"
$ws.Range("E37").Value = 46623.77

$ws.Range("A38").Value = 30597
$ws.Range("B38").Value = "Gina Craig"
$ws.Range("C38").Value = "elizabeth88@example.com"
$ws.Range("D38").Value = "Finance"
$ws.Range("E38").Value = 85393.8

$ws.Range("A39").Value = 3506
$ws.Range("B39").Value = "Devon Perez"
$ws.Range("C39").Value = "williamsellers@example.net"
$ws.Range("D39").Value = "IT"
$ws.Range("E39").Value = 74333.44

$ws.Range("A40").Value = 26501
$ws.Range("B40").Value = "Stephanie Shannon"
$ws.Range("C40").Value = "brianmassey@example.com"
$ws.Range("D40").Value = "Finance"
$ws.Range("E40").Value = 44799.56

$ws.Range("A41").Value = 9802
$ws.Range("B41").Value = "Stephanie Robinson"
$ws.Range("C41").Value = "wolfemichael@example.org"
$ws.Range("D41").Value = "Marketing"
$ws.Range("E41").Value = 87249.41

$ws.Range("A42").Value = 86072
$ws.Range("B42").Value = "Rachel Brown"
$ws.Range("C42").Value = "sarah15@example.org"
$ws.Range("D42").Value = "Finance"
$ws.Range("E42").Value = 85597.81

$ws.Range("A43").Value = 78407
$ws.Range("B43").Value = "David Oconnor"
$ws.Range("C43").Value = "bryanmorgan@example.net"
$ws.Range("D43").Value = "IT"
$ws.Range("E43").Value = 42880.66

$ws.Range("A44").Value = 85645
$ws.Range("B44").Value = "Eric Harris"
$ws.Range("C44").Value = "karen84@example.net"
$ws.Range("D44").Value = "IT"
$ws.Range("E44").Value = 56804.75

$ws.Range("A45").Value = 77525
$ws.Range("B45").Value = "Mrs. Taylor Morris"
$ws.Range("C45").Value = "taylorkimberly@example.net"
$ws.Range("D45").Value = "Marketing"
$ws.Range("E45").Value = 56055.9

$ws.Range("A46").Value = 30312
$ws.Range("B46").Value = "Deborah Lara"
$ws.Range("C46").Value = "mary84@example.net"
$ws.Range("D46").Value = "HR"
$ws.Range("E46").Value = 59789.99

$ws.Range("A47").Value = 89256
$ws.Range("B47").Value = "Kelly Ford"
$ws.Range("C47").Value = "rwatson@example.net"
$ws.Range("D47").Value = "HR"
$ws.Range("E47").Value = 88816.14

$ws.Range("A48").Value = 99618
$ws.Range("B48").Value = "Monique Edwards"
$ws.Range("C48").Value = "joshuagonzalez@example.com"
$ws.Range("D48").Value = "HR"
$ws.Range("E48").Value = 81739.34

$ws.Range("A49").Value = 64424
$ws.Range("B49").Value = "Daniel Douglas"
$ws.Range("C49").Value = "everettcarlos@example.net"
$ws.Range("D49").Value = "IT"
$ws.Range("E49").Value = 44695.78

$ws.Range("A50").Value = 84222
$ws.Range("B50").Value = "Jennifer Nelson"
$ws.Range("C50").Value = "millercaleb@example.net"
$ws.Range("D50").Value = "Finance"
$ws.Range("E50").Value = 97034.64999999999
